$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural change: insert new columns ---
# Insert 3 columns before old "contribuicoes" (M) to host apoio_std / apoio_min / apoio_max
$ws.Range("M1:O1").EntireColumn.Insert()
# Insert 3 columns after "contribuicoes_med" (now Q) to host contribuicoes_std/min/max
$ws.Range("R1:T1").EntireColumn.Insert()

# --- Rename headers (renamed metrics) ---
$ws.Range("H1").Value = "arrecadado_avg"
$ws.Range("I1").Value = "arrecadado_std"
$ws.Range("J1").Value = "arrecadado_min"
$ws.Range("K1").Value = "arrecadado_max"
$ws.Range("Q1").Value = "contribuicoes_med"

# --- New headers for the newly inserted columns ---
$ws.Range("M1").Value = "apoio_std"
$ws.Range("N1").Value = "apoio_min"
$ws.Range("O1").Value = "apoio_max"
$ws.Range("R1").Value = "contribuicoes_std"
$ws.Range("S1").Value = "contribuicoes_min"
$ws.Range("T1").Value = "contribuicoes_max"

# --- Row 2 (flex / apoia.se) new column values ---
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# --- Row 3 (flex / catarse): updated apoio_medio + new metric values ---
$ws.Range("L3").Value = 77.41063997458096
$ws.Range("M3").Value = 39.50983355883143
$ws.Range("N3").Value = 10.77163914429046
$ws.Range("O3").Value = 461.5197709071476

$ws.Range("R3").Value = 327.6748910926806
$ws.Range("S3").Value = 1
$ws.Range("T3").Value = 7954

# --- Number formats for the newly inserted columns (match neighbours) ---
$ws.Range("M2:O3").NumberFormat = "R$ #,##0.00"
$ws.Range("R2:T3").NumberFormat = "#,##0"
